# Regenerate the Handback status report: the handoff/handback round for the
# second tracked source file ("ae94a271-...") got re-run, producing a new
# pair of source-document GUIDs, a new content hash for the translated
# (.xlf) packages, and new handoff/handback timestamps. Propagate the new
# identifiers through the Overview sheet, the per-locale (zh-cn / de-de)
# detail sheets, and their hyperlinks' display text.

$wb = $excel.ActiveWorkbook

$oldGuid1 = "97f02eb6-3868-45bf-bdc9-eec5efc9cd25"      # row 2 source document id
$newGuid1 = "59db17ff-0d89-493f-a4fb-64bdf414a197"

$oldGuid2 = "ae94a271-8b44-4652-a391-beb04cb338c6"      # row 3 source document id
$newGuid2 = "ffffc5d5a3c0-02f9-4907-ab3a-fb673b2c366d"

$oldHash1 = "8a3a15b8aeeaa436431f53eb623dea5b0c7d03f7"  # row 2 .xlf package hash
$oldHash2 = "3db71f5cf1e4183c10e3cb06adb99d5c0c7b465e"  # row 3 .xlf package hash
$newHash  = "bac44bfe36bffe9cc476143af072f2ebdc47248a"  # new .xlf package hash (shared by both rows)

$zhHandoffOld = "2016-03-22 17:10:49"
$zhHandoffNew = "2016-03-22 17:12:55"
$zhHandbackOld = "2016-03-22 17:11:20"
$zhHandbackNew = "2016-03-22 17:13:23"

$deHandoffOld = "2016-03-22 17:10:54"
$deHandoffNew = "2016-03-22 17:12:59"
$deHandbackOld = "2016-03-22 17:11:27"
$deHandbackNew = "2016-03-22 17:13:30"

# --- 1) The "Correspond Handoff File" / "Correspond Handback File" (.xlf)
#        columns (D & G) on the row-3 record now reference the SAME source
#        GUID as row 2 (the two translation jobs were coalesced) -- swap the
#        row-3 GUID there to the row-2 GUID before the hash update so both
#        rows converge onto identical strings.
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("D3").Replace($oldGuid2, $newGuid1)
    $ws.Range("G3").Replace($oldGuid2, $newGuid1)
}

# --- 2) Update the .xlf content hash for both rows/columns so the two rows'
#        D/G values collapse to one shared string per locale sheet.
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("D2").Replace($oldHash1, $newHash)
    $ws.Range("G2").Replace($oldHash1, $newHash)
    $ws.Range("D3").Replace($oldHash2, $newHash)
    $ws.Range("G3").Replace($oldHash2, $newHash)
}

# --- 3) Update the handoff/handback timestamps (columns E & H) per locale.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2:E3").Replace($zhHandoffOld, $zhHandoffNew)
$wsZh.Range("H2:H3").Replace($zhHandbackOld, $zhHandbackNew)

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2:E3").Replace($deHandoffOld, $deHandoffNew)
$wsDe.Range("H2:H3").Replace($deHandbackOld, $deHandbackNew)

# --- 4) Update the source-document (.md) GUID everywhere it is still used as
#        a cell value: Overview!A, and the A/F columns on each locale sheet.
foreach ($sheetName in @("Overview", "zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Cells.Replace($oldGuid1, $newGuid1)
    $ws.Cells.Replace($oldGuid2, $newGuid2)
}

# --- 5) Hyperlink display text is a separate attribute from the cell value
#        and isn't touched by Range.Replace -- update each one explicitly.
#        ".xlf" links (columns D/G) always resolve to the row-2 GUID +
#        the new hash; ".md" links (columns A/F) keep each row's own GUID.
foreach ($sheetName in @("Overview", "zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($hl in $ws.Hyperlinks) {
        $t = $hl.TextToDisplay
        if ($t.Contains(".xlf")) {
            $t = $t.Replace($oldGuid1, $newGuid1)
            $t = $t.Replace($oldGuid2, $newGuid1)
            $t = $t.Replace($oldHash1, $newHash)
            $t = $t.Replace($oldHash2, $newHash)
        } else {
            $t = $t.Replace($oldGuid1, $newGuid1)
            $t = $t.Replace($oldGuid2, $newGuid2)
        }
        $hl.TextToDisplay = $t
    }
}
